$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in B6 and C6 while keeping their styling/borders
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# Update C7 value
$ws.Range("C7").Value = 33

# Update the active cell / selection to C6
$ws.Range("C6").Select()
